# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.550.43"
$ws.Range("E2").Value = "  +2.36%  "
$ws.Range("D3").Value = "3.649.11"
$ws.Range("E3").Value = "  +7.75%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "3.639.52"
$ws.Range("E7").Value = "  +7.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.622"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.32%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.610"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000285"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "4.237.32"
$ws.Range("E14").Value = "  +7.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "681.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.663.00"
$ws.Range("E17").Value = "  +8.13%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "71.708.92"
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.940"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +15.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "579.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("E35").Value = "  +2.02%  "
$ws.Range("E36").Value = "  +2.32%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.742.83"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.143"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  +3.61%  "
$ws.Range("E42").Value = "  +4.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0463"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.90%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.345"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.03%  "
$ws.Range("E48").Value = "  +3.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "133.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.44%  "
